$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The series in column A/B continues for one more quarter: a new last row
# (row 82) with the next date and its value. Copy the date formatting from
# the previous date cell (A81) so the new date cell keeps the workbook's
# existing custom date style instead of Excel creating a brand-new one.
$ws.Range("A81").Copy()
$ws.Range("A82").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row 82: date 2025-08-15 (serial 45884) and value 105.26
$ws.Cells.Item(82, 1).Value = 45884
$ws.Cells.Item(82, 2).Value = 105.26
